$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header / title strings (date 14.04.2025 -> 21.04.2025) ----
$ws.Range("A1").Value = "Mangrove Communication  21.04.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (21/04/2025) "

# ---- Top "Sale" table (rows 3-7) ----
# Row 3
$ws.Range("C3").Value = 30957
$ws.Range("D3").Value = 2826
$ws.Range("F3").Value = 2

# Row 4
$ws.Range("C4").Value = 23427
$ws.Range("D4").Value = 1413

# Row 5
$ws.Range("C5").Value = 23688
$ws.Range("D5").Value = 848
$ws.Range("F5").ClearContents()

# Row 6
$ws.Range("C6").Value = 31118
$ws.Range("D6").Value = 1413
$ws.Range("F6").Value = 2

# Row 7 (SUM formulas) recompute automatically

# ---- Stock table ----
# Row 14
$ws.Range("C14").Value = 217218
$ws.Range("D14").Value = 109190
$ws.Range("F14").ClearContents()

# Row 20
$ws.Range("C20").Value = 2000
$ws.Range("D20").Value = 2000

# Row 21
$ws.Range("C21").Value = 240
$ws.Range("D21").Value = 230

# Row 22
$ws.Range("C22").Value = 500
$ws.Range("D22").Value = 500

# Row 24
$ws.Range("C24").Value = 40

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 4

# Row 27
$ws.Range("C27").Value = 66
$ws.Range("D27").Value = 1

# Row 28
$ws.Range("C28").Value = 43

# ---- Bottom summary block ----
# Row 34
$ws.Range("H34").Value = 64678

# Row 35
$ws.Range("H35").Value = 17455

# Row 38
$ws.Range("H38").Value = 300000

# Insert a new row before the GRAND TOTAL row (old row 40), shifting GRAND
# TOTAL down to row 41, and give the new row the same formatting as the row
# above it (row 39).
$ws.Rows(40).Insert()
$ws.Range("F39:H39").Copy()
$ws.Range("F40:H40").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F40:G40").Merge()

$ws.Range("F40").Value = "Loan to MIST"
$ws.Range("H40").Value = 216000

# Extend the GRAND TOTAL formula (now on row 41) to include the new row 40
$ws.Range("H41").Formula = "=H33+H34+H35+H36+H37+H38+H39+H40"

# ---- Selection / scroll state ----
$ws.Range("H42").Select() | Out-Null
